# Fixed some bugs in AnalyzeReelsWaysEx3
# The test data rows got shuffled/corrected; update the A:F values for the
# affected rows so each row once again holds the correct symbol/reel counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(1001, 18, 30, 75, 60, 72)
    3  = @(601,  9,  60, 67, 60, 42)
    4  = @(1202, 2,  10, 10, 10, 10)
    5  = @(101,  9,  30, 15, 60, 15)
    6  = @(901,  16, 15, 45, 60, 60)
    7  = @(1201, 2,  10, 10, 10, 10)
    8  = @(902,  1,  0,  0,  0,  0)
    9  = @(201,  9,  30, 15, 45, 30)
    13 = @(501,  9,  52, 30, 75, 45)
    14 = @(401,  9,  48, 67, 75, 45)
    15 = @(701,  3,  90, 45, 97, 15)
    16 = @(2,    0,  2,  2,  2,  2)
    17 = @(3,    0,  3,  3,  3,  3)
    19 = @(502,  0,  4,  0,  0,  0)
    20 = @(1101, 0,  15, 30, 30, 0)
    21 = @(802,  0,  4,  5,  4,  0)
}

$cols = @("A", "B", "C", "D", "E", "F")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $values[$i]
    }
}
